$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Workbook window size (bookViews/workbookView windowWidth/windowHeight)
# Best-effort; the underlying engine may not persist this, but set it anyway.
# ---------------------------------------------------------------------------
try {
    $win = $excel.ActiveWindow
    $win.Width = 12075
    $win.Height = 9600
} catch {}

# ---------------------------------------------------------------------------
# Row 31: clear "Spring" / "PullSpring, PushSpring" (this pair moves to row 35)
# ---------------------------------------------------------------------------
$ws.Range("A31").ClearContents()
$ws.Range("B31").ClearContents()

# ---------------------------------------------------------------------------
# Row 32: drop the stray "Constant" label in D32
# ---------------------------------------------------------------------------
$ws.Range("D32").ClearContents()

# ---------------------------------------------------------------------------
# Row 36: this whole row (Gravity / GravityLabel / Connect... afield) is removed;
# its content is replaced by new rows 37/38 below, with an updated sentence.
# ---------------------------------------------------------------------------
$ws.Range("A36").ClearContents()
$ws.Range("D36").ClearContents()
$ws.Range("F36").ClearContents()

# ---------------------------------------------------------------------------
# Row 35 (new): "Spring" / "PullSpring, PushSpring" -- same look as old row 31
# ---------------------------------------------------------------------------
$ws.Range("A35").Value = "Spring"
$ws.Range("A35").Font.Color = 192
$ws.Range("B35").Value = "PullSpring, PushSpring"

# ---------------------------------------------------------------------------
# Row 37 (new): Gravity / GravityLabel / updated sentence ("a field")
# ---------------------------------------------------------------------------
$ws.Range("A37").Value = "Gravity"
$ws.Range("D37").Value = "GravityLabel"
$ws.Range("F37").Value = "Connect this node to the gravitywell with GravityLabel with a field"

# ---------------------------------------------------------------------------
# Row 38 (new): GravityWell
# ---------------------------------------------------------------------------
$ws.Range("A38").Value = "GravityWell"

# ---------------------------------------------------------------------------
# New rows 40, 41, 43, 44
# ---------------------------------------------------------------------------
$ws.Range("A40").Value = "Stationary"
$ws.Range("A41").Value = "NodeAngle"
$ws.Range("A43").Value = "Visual"
$ws.Range("A44").Value = "Colour"

# ---------------------------------------------------------------------------
# Sheet view: scroll so row 4 is at top, and select A45
# ---------------------------------------------------------------------------
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 4
    $win.ScrollColumn = 1
} catch {}
$ws.Range("A45").Select()
